# Append a new Adafruit IO feed reading as row 60 of the sheet, mirroring
# the existing "temperature" rows (timestamp, feed key, value, lat/long/elev).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 60

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# The "Value" column holds numeric-looking readings as plain text in this
# sheet (e.g. "25"), same as every prior row. Force text so Excel doesn't
# silently coerce it to a number, then drop the temporary text format so
# the cell is left without any explicit style, matching the rest of the
# sheet.
$valueCell = $ws.Cells.Item($row, 3)
$valueCell.NumberFormat = "@"
$valueCell.Value = "25"
$valueCell.ClearFormats()

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
